$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1200
$ws.Range("J97").Value = 1200
$ws.Range("L97").Value = 3600
$ws.Range("N97").Value = -4592

$ws.Range("H100").Value = 5490
$ws.Range("I100").Value = 1000
$ws.Range("J100").Value = 9980
$ws.Range("K100").Value = 1000
$ws.Range("L100").Value = 9980
$ws.Range("M100").Value = -459
$ws.Range("N100").Value = -11062

$ws.Range("H112").Value = 7554.3335
$ws.Range("J112").Value = 7554.3335
$ws.Range("L112").Value = 22663.0005
$ws.Range("N112").Value = -24879.0005

$ws.Range("H116").Value = 7837.4546
$ws.Range("I116").Value = 13263.223
$ws.Range("K116").Value = 13263.223
$ws.Range("M116").Value = -9821.223

$ws.Range("H137").Value = 2127.4
$ws.Range("I137").Value = 1726.2222
$ws.Range("J137").Value = 2455.6365
$ws.Range("K137").Value = 5178.6666
$ws.Range("L137").Value = 7366.9095
$ws.Range("M137").Value = -2628.6666
$ws.Range("N137").Value = -12466.9095

$ws.Range("H138").Value = 6954.5835
$ws.Range("I138").Value = 6765.778
$ws.Range("J138").Value = 7017.5186
$ws.Range("K138").Value = 20297.334
$ws.Range("L138").Value = 21052.5558
$ws.Range("M138").Value = -15157.334
$ws.Range("N138").Value = -31332.5558

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2994.9194
$ws.Range("I32").Value = 2406.0266
$ws.Range("K32").Value = 2406.0266
$ws.Range("M32").Value = -2119.0266

$ws.Range("H61").Value = 3553.7778
$ws.Range("I61").Value = 2120.1
$ws.Range("K61").Value = 2120.1
$ws.Range("M61").Value = -1908.1

$ws.Range("H132").Value = 3876.1875
$ws.Range("I132").Value = 2940.125
$ws.Range("J132").Value = 4812.25
$ws.Range("K132").Value = 8820.375
$ws.Range("L132").Value = 14436.75
$ws.Range("M132").Value = -6290.375
$ws.Range("N132").Value = -19496.75

$ws.Range("H136").Value = 3553.7778
$ws.Range("I136").Value = 2120.1
$ws.Range("K136").Value = 6360.299999999999
$ws.Range("M136").Value = -3810.299999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 6866.8
$ws.Range("I80").Value = 27
$ws.Range("J80").Value = 7919.077
$ws.Range("K80").Value = 27
$ws.Range("L80").Value = 7919.077
$ws.Range("M80").Value = 971
$ws.Range("N80").Value = -9915.077000000001

$ws.Range("H83").Value = 6866.8
$ws.Range("I83").Value = 27
$ws.Range("J83").Value = 7919.077
$ws.Range("K83").Value = 135
$ws.Range("L83").Value = 39595.385
$ws.Range("M83").Value = 4857
$ws.Range("N83").Value = -49579.385

$ws.Range("H102").Value = 11900
$ws.Range("I102").Value = 1850
$ws.Range("J102").Value = 32000
$ws.Range("K102").Value = 1850
$ws.Range("L102").Value = 32000
$ws.Range("M102").Value = 1395
$ws.Range("N102").Value = -38490

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1194.1758
$ws.Range("I31").Value = 836.75
$ws.Range("J31").Value = 1474.5098
$ws.Range("K31").Value = 836.75
$ws.Range("L31").Value = 1474.5098
$ws.Range("M31").Value = -541.75
$ws.Range("N31").Value = -2064.5098

$ws.Range("H34").Value = 1194.1758
$ws.Range("I34").Value = 836.75
$ws.Range("J34").Value = 1474.5098
$ws.Range("K34").Value = 836.75
$ws.Range("L34").Value = 1474.5098
$ws.Range("M34").Value = -634.75
$ws.Range("N34").Value = -1878.5098

$ws.Range("H99").Value = 1982
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 1976
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 1976
$ws.Range("M99").Value = -502
$ws.Range("N99").Value = -4972

$ws.Range("H126").Value = 1982
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 1976
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 5928
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -10868

$ws.Range("H130").Value = 23800
$ws.Range("J130").Value = 23800
$ws.Range("L130").Value = 23800
$ws.Range("N130").Value = -33840

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3224.9246
$ws.Range("I68").Value = 1345.5333
$ws.Range("J68").Value = 3966.7896
$ws.Range("K68").Value = 4036.5999
$ws.Range("L68").Value = 11900.3688
$ws.Range("M68").Value = -3225.5999
$ws.Range("N68").Value = -13522.3688

$ws.Range("H71").Value = 3224.9246
$ws.Range("I71").Value = 1345.5333
$ws.Range("J71").Value = 3966.7896
$ws.Range("K71").Value = 12109.7997
$ws.Range("L71").Value = 35701.1064
$ws.Range("M71").Value = -8053.7997
$ws.Range("N71").Value = -43813.1064

$ws.Range("H107").Value = 4095.4736
$ws.Range("J107").Value = 4156.3335
$ws.Range("L107").Value = 12469.0005
$ws.Range("N107").Value = -16309.0005

$ws.Range("H132").Value = 1884.8572
$ws.Range("J132").Value = 2623.5
$ws.Range("L132").Value = 23611.5
$ws.Range("N132").Value = -28671.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 206
$ws.Range("I2").Value = 258.75
$ws.Range("J2").Value = 163.8
$ws.Range("K2").Value = 258.75
$ws.Range("L2").Value = 163.8
$ws.Range("M2").Value = -145.75
$ws.Range("N2").Value = -389.8

$ws.Range("H132").Value = 1541798.6
$ws.Range("I132").Value = 2139061.8
$ws.Range("J132").Value = 5979.4287
$ws.Range("K132").Value = 6417185.399999999
$ws.Range("L132").Value = 17938.2861
$ws.Range("M132").Value = -6414655.399999999
$ws.Range("N132").Value = -22998.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8217
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 8217
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 8217
$ws.Range("N22").Value = -8807
$ws.Range("M22").ClearContents()

$ws.Range("H27").Value = 8217
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 8217
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 8217
$ws.Range("N27").Value = -8431
$ws.Range("M27").ClearContents()

$ws.Range("H94").Value = 44859.668
$ws.Range("J94").Value = 44859.668
$ws.Range("L94").Value = 44859.668
$ws.Range("N94").Value = -46211.668

$ws.Range("H100").Value = 1463.1538
$ws.Range("I100").Value = 1185.5454
$ws.Range("J100").Value = 2990
$ws.Range("K100").Value = 1185.5454
$ws.Range("L100").Value = 2990
$ws.Range("M100").Value = -644.5454
$ws.Range("N100").Value = -4072

$ws.Range("H134").Value = 43613.6
$ws.Range("J134").Value = 43613.6
$ws.Range("L134").Value = 43613.6
$ws.Range("N134").Value = -53753.6

$ws.Range("H136").Value = 4043.5
$ws.Range("I136").Value = 2978.9666
$ws.Range("J136").Value = 6704.8335
$ws.Range("K136").Value = 8936.899800000001
$ws.Range("L136").Value = 20114.5005
$ws.Range("M136").Value = -6386.899800000001
$ws.Range("N136").Value = -25214.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 9579.362999999999
$ws.Range("I126").Value = 9985.944
$ws.Range("J126").Value = 7749.75
$ws.Range("K126").Value = 29957.832
$ws.Range("L126").Value = 23249.25
$ws.Range("M126").Value = -27487.832
$ws.Range("N126").Value = -28189.25
